$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 3
$ws.Range("A3").Value = "2022-09-20 13:35:17"
$ws.Range("B3").Value = "MER_SMI_1Junho_8"
$ws.Range("C3").Value = "MER SMI"

# Update row 4
$ws.Range("A4").Value = "2022-09-20 13:35:18"
$ws.Range("B4").Value = "MER_SMI_1Junho_8"
$ws.Range("C4").Value = "MER SMI"

# Update row 5
$ws.Range("A5").Value = "2022-09-20 13:35:18"
$ws.Range("B5").Value = "MER_SMI_1Junho_8"
$ws.Range("C5").Value = "MER SMI"
$ws.Range("D5").Value = "Buscar valores para cada indicador: DSD PMTCT EID"

# Update row 6
$ws.Range("A6").Value = "2022-09-20 13:35:19"
$ws.Range("B6").Value = "MER_SMI_1Junho_8"
$ws.Range("C6").Value = "MER SMI"
$ws.Range("D6").Value = "Buscar valores para cada indicador: DSD PMTCT HEI POS"

# Remove rows 7 through 11 (previously existing extra rows)
$ws.Range("A7:E11").Delete()
